$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.253.62"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "2.847.14"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "361.40"
$ws.Range("E5").Value = "  +5.87%  "
$ws.Range("D6").Value = "113.64"
$ws.Range("E6").Value = "  -3.14%  "
$ws.Range("D7").Value = "0.574"
$ws.Range("E7").Value = "  +3.79%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "0.604"
$ws.Range("E9").Value = "  +3.87%  "
$ws.Range("D10").Value = "41.64"
$ws.Range("E10").Value = "  -1.41%  "
$ws.Range("D11").Value = "0.0863"
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("D12").Value = "0.131"
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("D13").Value = "20.00"
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("E14").Value = "  +1.72%  "
$ws.Range("D15").Value = "3.295.03"
$ws.Range("E15").Value = "  +1.82%  "
$ws.Range("D16").Value = "2.858.22"
$ws.Range("E16").Value = "  +2.73%  "
$ws.Range("D17").Value = "0.904"
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("D18").Value = "52.182.55"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").Value = "7.59"
$ws.Range("E19").Value = "  +8.91%  "
$ws.Range("D20").Value = "3.16"
$ws.Range("E20").Value = "  -2.35%  "
$ws.Range("D21").Value = "13.62"
$ws.Range("E21").Value = "  +1.38%  "
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").Value = "70.34"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "268.39"
$ws.Range("E24").Value = "  -3.97%  "
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("D26").Value = "27.24"
$ws.Range("E26").Value = "  +0.86%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").Value = "10.44"
$ws.Range("E28").Value = "  +1.92%  "
$ws.Range("E29").Value = "  +1.47%  "
$ws.Range("D30").Value = "53.88"
$ws.Range("E30").Value = "  +6.83%  "
$ws.Range("E31").Value = "  -1.50%  "
$ws.Range("D32").Value = "34.44"
$ws.Range("E32").Value = "  -1.66%  "
$ws.Range("D33").Value = "0.0459"
$ws.Range("E33").Value = "  +23.29%  "
$ws.Range("E34").Value = "  +2.52%  "
$ws.Range("D35").Value = "5.41"
$ws.Range("E35").Value = "  +8.09%  "
$ws.Range("E36").Value = "  +2.01%  "
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").Value = "3.28"
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("E39").Value = "  -2.69%  "
$ws.Range("D40").Value = "18.39"
$ws.Range("E40").Value = "  -3.70%  "
$ws.Range("D41").Value = "23.86"
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("D42").Value = "0.117"
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("D43").Value = "128.61"
$ws.Range("E43").Value = "  +1.71%  "
$ws.Range("E44").Value = "  -6.75%  "
$ws.Range("E45").Value = "  -2.76%  "
$ws.Range("E46").Value = "  +1.41%  "
$ws.Range("D47").Value = "2.119.16"
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("E49").Value = "  +9.48%  "
$ws.Range("E50").Value = "  +5.27%  "
$ws.Range("E51").Value = "  +1.31%  "
